$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2760408939892198
$ws.Range("C2").Value = 0.2760408939892198

$ws.Range("B3").Value = 0.5424796504655186
$ws.Range("C3").Value = 0.5424796504655186

$ws.Range("B4").Value = 0.7170038108021736
$ws.Range("C4").Value = 0.7170038108021736

$ws.Range("B5").Value = 0.009443920072667887
$ws.Range("C5").Value = 0.004353953684610743
$ws.Range("D5").Value = 241
$ws.Range("E5").Value = 192

$ws.Range("B6").Value = 0.9391719980261437
$ws.Range("C6").Value = 0.9391719980261437

$ws.Range("B7").Value = 0.00005322173314725846
$ws.Range("C7").Value = 0.00004552711771745896
$ws.Range("D7").Value = 825
$ws.Range("E7").Value = 797

$ws.Range("B8").Value = 0.1011464052136305
$ws.Range("C8").Value = 0.07995000285879053
$ws.Range("D8").Value = 146
$ws.Range("E8").Value = 123

$ws.Range("B9").Value = 0.03507900269087773
$ws.Range("C9").Value = 0.0329916560106028
$ws.Range("D9").Value = 282
$ws.Range("E9").Value = 271

$ws.Range("B10").Value = 0.3970927852420561
$ws.Range("C10").Value = 0.3857449803705244
$ws.Range("D10").Value = 37
$ws.Range("G10").Value = 124

$ws.Range("B11").Value = 0.000000004020466781811809
$ws.Range("C11").Value = 0.000000002528190560195616
$ws.Range("D11").Value = 1099
$ws.Range("E11").Value = 1051
$ws.Range("F11").Value = 1484
$ws.Range("G11").Value = 1481

$ws.Range("B12").Value = 0.07952597316640113
$ws.Range("C12").Value = 0.07952597316640113

$ws.Range("B13").Value = 0.507471706605126
$ws.Range("C13").Value = 0.5004452583181059
$ws.Range("D13").Value = 90
$ws.Range("E13").Value = 81

$ws.Range("B14").Value = 0.7688953153587761
$ws.Range("C14").Value = 0.7617814239549737
$ws.Range("D14").Value = 30

$ws.Range("B15").Value = 0.7894068617535623
$ws.Range("C15").Value = 0.7868686898597245
$ws.Range("D15").Value = 18
$ws.Range("G15").Value = 31

$ws.Range("B16").Value = 0.6266129947137067
$ws.Range("C16").Value = 0.6266129947137067

$ws.Range("B17").Value = 0.6448527714974646
$ws.Range("C17").Value = 0.6448527714974646
